$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data to add to the Management section (rows 8-12)
$data = @(
    @("Member don’t care about Risks in project and don’t update Risk Category", "Follow risk plan, plan about interation for risk", "Knowing more about manage risk better", "Management"),
    @("Too much and more difficult to measurement", "Research about measurement, implement Goal-Question-Metric", "Knowing more defenite about metrics and how to get it", "Management"),
    @("Project difficult to control and monitoring", "Plan for detail plan, WBS, implement tracking and monitoring through measurement about schedule devition metric", "", "Management"),
    @("360 review is not good conduct", "Require team member write reflection base on 360 reivew", "Knowing about management and communicate between team member", "Management"),
    @("Team member is not complete work on time", "Re-estimate, and evaluate effort of team member", "Conduct measurement about productivity", "Management")
)

$row = 8
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = ($row - 2)
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    if ($entry[2] -ne "") {
        $ws.Cells.Item($row, 4).Value = $entry[2]
    }
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row++
}

$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B10").Select()
